$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.302.18"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.913.11"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "0.725"
$ws.Range("E5").Value = "  +9.41%  "
$ws.Range("D6").Value = "255.66"
$ws.Range("E6").Value = "  +3.81%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").Value = "42.17"
$ws.Range("E8").Value = "  +0.98%  "
$ws.Range("D9").Value = "0.369"
$ws.Range("E9").Value = "  +6.32%  "
$ws.Range("D10").Value = "53.32"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("E11").Value = "  +6.42%  "
$ws.Range("D12").Value = "0.0990"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "13.18"
$ws.Range("E13").Value = "  +7.22%  "
$ws.Range("D14").Value = "2.191.39"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "0.738"
$ws.Range("E15").Value = "  +5.36%  "
$ws.Range("D16").Value = "5.01"
$ws.Range("E16").Value = "  +3.77%  "
$ws.Range("D17").Value = "1.939.04"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "35.300.20"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "75.25"
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("D20").Value = "0.0₃0850"
$ws.Range("E20").Value = "  +3.73%  "
$ws.Range("D21").Value = "246.08"
$ws.Range("E21").Value = "  +1.96%  "
$ws.Range("D22").Value = "13.12"
$ws.Range("E22").Value = "  +5.05%  "
$ws.Range("E23").Value = "  +6.74%  "
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("E25").Value = "  +8.09%  "
$ws.Range("D26").Value = "2.40"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").Value = "166.68"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("D28").Value = "8.81"
$ws.Range("E28").Value = "  +4.35%  "
$ws.Range("D29").Value = "18.82"
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("E30").Value = "  +4.45%  "
$ws.Range("D31").Value = "4.127.51"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").Value = "1.70"
$ws.Range("E32").Value = "  +27.35%  "
$ws.Range("E33").Value = "  +5.56%  "
$ws.Range("E34").Value = "  +14.53%  "
$ws.Range("E35").Value = "  +4.66%  "
$ws.Range("E36").Value = "  +4.23%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").Value = "0.920"
$ws.Range("E38").Value = "  -2.82%  "
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "100.06"
$ws.Range("E40").Value = "  +11.23%  "
$ws.Range("D41").Value = "0.0222"
$ws.Range("E41").Value = "  +6.20%  "
$ws.Range("D42").Value = "17.07"
$ws.Range("E42").Value = "  +4.96%  "
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").Value = "1.344.67"
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("D48").Value = "6.75"
$ws.Range("E48").Value = "  +3.53%  "
$ws.Range("E49").Value = "  -1.08%  "
$ws.Range("D50").Value = "45.09"
$ws.Range("E50").Value = "  -8.47%  "
$ws.Range("D51").Value = "0.0758"
$ws.Range("E51").Value = "  +7.02%  "